# Regenerate save_data column G ("K") values to reflect the new
# Strike-count-based ("K") calculation instead of the old Strike# values.
# The mapping below gives the new value for each data row (row 2 = first
# data row, i.e. index 0) in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2 = 1; 3 = 2; 4 = 0; 5 = 0; 6 = 1; 7 = 1; 8 = 2; 9 = 1; 10 = 1;
    11 = 1; 12 = 0; 13 = 2; 14 = 0; 15 = 0; 16 = 1; 17 = 0; 18 = 0;
    19 = 3; 20 = 0; 21 = 1; 23 = 1; 24 = 0; 25 = 3; 27 = 1; 28 = 2;
    29 = 1; 30 = 3; 31 = 1; 32 = 1; 33 = 0; 34 = 1; 35 = 1; 36 = 2;
    37 = 2; 38 = 0; 39 = 3; 40 = 0; 41 = 0; 42 = 2; 43 = 3; 44 = 1;
    45 = 1; 46 = 1; 47 = 2; 48 = 0; 49 = 1; 50 = 0; 51 = 1; 52 = 2;
    53 = 1; 54 = 1; 55 = 2; 57 = 0; 58 = 2; 59 = 2; 60 = 1; 61 = 1;
    62 = 2; 63 = 1; 64 = 1; 66 = 1; 67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
